$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.041452
$ws.Cells.Item(2, 8).Value = 0.124356
$ws.Cells.Item(2, 9).Value = 0.05439747478414846
$ws.Cells.Item(2, 10).Value = 0.05439747478414846
$ws.Cells.Item(2, 13).Value = 0.09551033333333332
$ws.Cells.Item(2, 14).Value = 0.286531
$ws.Cells.Item(2, 15).Value = 0.0198020999427218
$ws.Cells.Item(2, 16).Value = 0.0198020999427218
$ws.Cells.Item(2, 17).Value = 0.003959094337333333
$ws.Cells.Item(2, 18).Value = 0.035631849036
$ws.Cells.Item(2, 19).Value = 0.001077184232307397
$ws.Cells.Item(2, 20).Value = 0.001077184232307397
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.041452
$ws.Cells.Item(3, 8).Value = 0.124356
$ws.Cells.Item(3, 9).Value = 0.05439747478414846
$ws.Cells.Item(3, 10).Value = 0.05439747478414846
$ws.Cells.Item(3, 15).Value = 0.07175622098770619
$ws.Cells.Item(3, 16).Value = 0.07175622098770619
$ws.Cells.Item(3, 17).Value = 0.01434644047866667
$ws.Cells.Item(3, 18).Value = 0.129117964308
$ws.Cells.Item(3, 19).Value = 0.003903357221784532
$ws.Cells.Item(3, 20).Value = 0.003903357221784532
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.041452
$ws.Cells.Item(4, 8).Value = 0.124356
$ws.Cells.Item(4, 9).Value = 0.05439747478414846
$ws.Cells.Item(4, 10).Value = 0.05439747478414846
$ws.Cells.Item(4, 13).Value = 4.381634666666667
$ws.Cells.Item(4, 14).Value = 13.144904
$ws.Cells.Item(4, 15).Value = 0.9084416790695721
$ws.Cells.Item(4, 16).Value = 0.9084416790695721
$ws.Cells.Item(4, 17).Value = 0.1816275202026667
$ws.Cells.Item(4, 18).Value = 1.634647681824
$ws.Cells.Item(4, 19).Value = 0.04941693333005653
$ws.Cells.Item(4, 20).Value = 0.04941693333005653
$ws.Cells.Item(5, 9).Value = 0.5204718857143857
$ws.Cells.Item(5, 10).Value = 0.5204718857143856
$ws.Cells.Item(5, 13).Value = 0.09551033333333332
$ws.Cells.Item(5, 14).Value = 0.286531
$ws.Cells.Item(5, 15).Value = 0.0198020999427218
$ws.Cells.Item(5, 16).Value = 0.0198020999427218
$ws.Cells.Item(5, 17).Value = 0.03788038514011111
$ws.Cells.Item(5, 18).Value = 0.340923466261
$ws.Cells.Item(5, 19).Value = 0.01030643629829314
$ws.Cells.Item(5, 20).Value = 0.01030643629829314
$ws.Cells.Item(6, 9).Value = 0.5204718857143857
$ws.Cells.Item(6, 10).Value = 0.5204718857143856
$ws.Cells.Item(6, 15).Value = 0.07175622098770619
$ws.Cells.Item(6, 16).Value = 0.07175622098770619
$ws.Cells.Item(6, 19).Value = 0.03734709564920962
$ws.Cells.Item(6, 20).Value = 0.03734709564920961
$ws.Cells.Item(7, 9).Value = 0.5204718857143857
$ws.Cells.Item(7, 10).Value = 0.5204718857143856
$ws.Cells.Item(7, 13).Value = 4.381634666666667
$ws.Cells.Item(7, 14).Value = 13.144904
$ws.Cells.Item(7, 15).Value = 0.9084416790695721
$ws.Cells.Item(7, 16).Value = 0.9084416790695721
$ws.Cells.Item(7, 17).Value = 1.737801585691556
$ws.Cells.Item(7, 18).Value = 15.640214271224
$ws.Cells.Item(7, 19).Value = 0.472818353766883
$ws.Cells.Item(7, 20).Value = 0.4728183537668829
$ws.Cells.Item(8, 7).Value = 0.3239583333333333
$ws.Cells.Item(8, 8).Value = 0.971875
$ws.Cells.Item(8, 9).Value = 0.4251306395014658
$ws.Cells.Item(8, 10).Value = 0.4251306395014658
$ws.Cells.Item(8, 13).Value = 0.09551033333333332
$ws.Cells.Item(8, 14).Value = 0.286531
$ws.Cells.Item(8, 15).Value = 0.0198020999427218
$ws.Cells.Item(8, 16).Value = 0.0198020999427218
$ws.Cells.Item(8, 17).Value = 0.03094136840277778
$ws.Cells.Item(8, 18).Value = 0.278472315625
$ws.Cells.Item(8, 19).Value = 0.008418479412121257
$ws.Cells.Item(8, 20).Value = 0.008418479412121257
$ws.Cells.Item(9, 7).Value = 0.3239583333333333
$ws.Cells.Item(9, 8).Value = 0.971875
$ws.Cells.Item(9, 9).Value = 0.4251306395014658
$ws.Cells.Item(9, 10).Value = 0.4251306395014658
$ws.Cells.Item(9, 15).Value = 0.07175622098770619
$ws.Cells.Item(9, 16).Value = 0.07175622098770619
$ws.Cells.Item(9, 17).Value = 0.1121212232638889
$ws.Cells.Item(9, 18).Value = 1.009091009375
$ws.Cells.Item(9, 19).Value = 0.03050576811671204
$ws.Cells.Item(9, 20).Value = 0.03050576811671203
$ws.Cells.Item(10, 7).Value = 0.3239583333333333
$ws.Cells.Item(10, 8).Value = 0.971875
$ws.Cells.Item(10, 9).Value = 0.4251306395014658
$ws.Cells.Item(10, 10).Value = 0.4251306395014658
$ws.Cells.Item(10, 13).Value = 4.381634666666667
$ws.Cells.Item(10, 14).Value = 13.144904
$ws.Cells.Item(10, 15).Value = 0.9084416790695721
$ws.Cells.Item(10, 16).Value = 0.9084416790695721
$ws.Cells.Item(10, 17).Value = 1.419467063888889
$ws.Cells.Item(10, 18).Value = 12.775203575
$ws.Cells.Item(10, 19).Value = 0.3862063919726326
$ws.Cells.Item(10, 20).Value = 0.3862063919726325
